$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cell D9 currently holds a lone space " " (a placeholder shared string).
# Re-purpose that same text slot by overwriting it with "middleware" first,
# then copy the resulting text across to A8 - this way the underlying shared
# string entry is edited in place (rather than appending a brand new one).
$ws.Range("D9").Value = "middleware"
$ws.Range("A8").Value = $ws.Range("D9").Text

# Fill in the rest of the new "middleware" row.
$ws.Range("B8").Value = "opogramowanie pośredniczące"
$ws.Range("C8").Value = "populate"
$ws.Range("D8").Value = "zaludniać/zasiedlać"

# New row 9: colon / dwukropek / be aware of / ...
$ws.Range("A9").Value = "colon"
$ws.Range("B9").Value = "dwukropek"
$ws.Range("C9").Value = "be aware of"

# New row 10: crucial / kluczowy (only columns C and D populated).
$ws.Range("C10").Value = "crucial"
$ws.Range("D10").Value = "kluczowy"

# Finally, restore D9 to its real value (reuses the already-existing
# "co ciekawe" shared string).
$ws.Range("D9").Value = "co ciekawe"

# Match the saved selection/active cell from the edited workbook.
$ws.Range("D10").Select()
